# Update the dSF ("F") column values for Zack Greinke's 2021 game log.
# These reflect a data repull where dSF (change in win probability / score
# differential from final state) differs from dS0 for several starts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -1
    4  = -3
    6  = -2
    7  = -5
    8  = -3
    9  = -1
    10 = -1
    11 = 4
    13 = -1
    14 = 2
    15 = 3
    16 = -1
    18 = -1
    19 = 2
    20 = 4
    21 = -1
    22 = 12
    23 = 6
    24 = -4
    25 = 7
    26 = 7
    27 = 4
    28 = 0
    29 = -2
    30 = -2
    31 = 1
    32 = -5
    34 = 3
    35 = -2
    38 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
